# correção nos dados e inicio da analise PNAD 2009
# Fix header labels in row 2 that were previously unnamed/placeholder
# (pandas auto-generated "unnamed: N_level_1" labels) and should read "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
